$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-08-23"

# Update the label for the August row
$ws.Range("A9").Value = "August (through 08-23)"

# Update the August (row 9) values for columns C..I (2016..2022)
$ws.Range("C9").Value = 52
$ws.Range("D9").Value = 62
$ws.Range("E9").Value = 41
$ws.Range("F9").Value = 34
$ws.Range("G9").Value = 135
$ws.Range("H9").Value = 123
$ws.Range("I9").Value = 130

# Update the Total (row 10) values for columns C..I (2016..2022)
$ws.Range("C10").Value = 354
$ws.Range("D10").Value = 527
$ws.Range("E10").Value = 466
$ws.Range("F10").Value = 338
$ws.Range("G10").Value = 756
$ws.Range("H10").Value = 1033
$ws.Range("I10").Value = 1101
